$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all affected cells (rows 3-21, columns A-N) to Text format to preserve numeric-looking strings as text
$ws.Range("A3:N21").NumberFormat = "@"

# Row 3 <= original row 6 (url "https://www.immowelt.de/expose/2an6q5g")
$ws.Cells.Item(3, 1).Value = '"https://www.immowelt.de/expose/2an6q5g"'
$ws.Cells.Item(3, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(3, 3).Value = '"Wohnberechtigungsschein (WBS3) erforderlich!! Gem\u00fctliche Zweiraumwohnung auf 59m\u00b2"'
$ws.Cells.Item(3, 4).Value = '566.96'
$ws.Cells.Item(3, 5).Value = '"EUR"'
$ws.Cells.Item(3, 6).Value = '59.68'
$ws.Cells.Item(3, 7).Value = '2'
$ws.Cells.Item(3, 8).Value = '2019'
$ws.Cells.Item(3, 9).Value = '["Stellplatz"]'
$ws.Cells.Item(3, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 15", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(3, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1700", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 566.96, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 199, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 124, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 949.96, "Unit": "EUR"}]}'
$ws.Cells.Item(3, 12).Value = '["https://ms.immowelt.org/04bb8e8f-62a6-4aa1-be23-fcac75a2ca80/64e8eae5-448b-4e0a-aef6-3ee4487d6cad", "https://ms.immowelt.org/07409987-f4b5-4a2f-9ed5-55d000afa646/e6d2fc11-39c9-44c4-9ab7-6b79a5e4667d", "https://ms.immowelt.org/d8d286af-d70c-4388-a156-0e4525e576c2/3f08e4d7-3ff5-4332-b507-dfc163a7d33f", "https://ms.immowelt.org/0eb043d6-e29a-4797-8d47-b95fdf183238/bda6c2c7-7309-497f-836d-6483f70abd43", "https://ms.immowelt.org/4aed45df-9c7d-4a2b-b188-b8968b22727d/50fd789e-2343-4b28-867b-f726db140de4"]'
$ws.Cells.Item(3, 13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(3, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnberechtigungsschein (WBS3) erforderlich!! Gem\u00fctliche Zweiraumwohnung auf 59m\u00b2", "basicObjectPricEur": 566.96, "basicLivingSpace": 59.68, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 15", "basicUrl": "https://www.immowelt.de/expose/2an6q5g"}'

# Row 4 <= original row 8 (url "https://www.immowelt.de/expose/25dtq5g")
$ws.Cells.Item(4, 1).Value = '"https://www.immowelt.de/expose/25dtq5g"'
$ws.Cells.Item(4, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(4, 3).Value = '"3-ZW W\u00fcrzburg  Sanderau,  3. OG.  Bad neu."'
$ws.Cells.Item(4, 4).Value = '950'
$ws.Cells.Item(4, 5).Value = '"EUR"'
$ws.Cells.Item(4, 6).Value = '79'
$ws.Cells.Item(4, 7).Value = '3'
$ws.Cells.Item(4, 8).Value = '1911'
$ws.Cells.Item(4, 9).Value = '["Altbau (bis 1945)", "WANNE", "FENSTER", "Kelleranteil", "GAS", "Etagenheizung", "Einbauk\u00fcche", "Balkon"]'
$ws.Cells.Item(4, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Sanderau", "ZipCode": "97072", "Street": "Eichendorffstrasse 4", "LocationId": 496022, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(4, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1900", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Sanderau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg/mietspiegel"}}, "DataTable": [{"NumberValue": 950, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 120, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}]}'
$ws.Cells.Item(4, 12).Value = '["https://ms.immowelt.org/53a40e8a-24fb-45ae-a531-fff2a2f4ef38/73c3faa7-d432-44a3-995c-dd6a57b55ff8", "https://ms.immowelt.org/5fd3d699-98a8-4599-a121-418d851bdca3/552eef5e-2a9f-4c01-b006-8b517f701ec2", "https://ms.immowelt.org/e1c22ef7-fa4f-4bbe-b945-4a0c9417515d/e28a4ad5-ff43-4cdd-a230-41a8a963e3de", "https://ms.immowelt.org/1fe9ddd3-4952-4f94-9aba-7b98e8a565e2/4cc8d480-3a53-414c-ada9-ef0da80900ed", "https://ms.immowelt.org/2e8ae348-cab2-4e79-8266-ee52099c1e73/524bcb38-9f3e-44ef-afa6-8810bc8ec8bf", "https://ms.immowelt.org/b0e64c11-6bbf-4ed6-9eb5-3767349429ee/e3628c5a-5b53-4019-a135-ff500cf7a0ac", "https://ms.immowelt.org/4a8c3d49-12ed-417a-b206-b38a38d9e894/7e1402ca-1074-4126-8ddc-ba9dbe8eac83", "https://ms.immowelt.org/82deeace-1bd0-4723-8a33-2c2ad6084c4a/cabe6e07-e983-4c08-826a-e189aa3d4288", "https://ms.immowelt.org/4e16951a-d91d-485c-a019-726374d3e5d3/97f13c47-29c0-4a36-8ee3-73268762ea32", "https://ms.immowelt.org/f8cc63cd-6dc9-4cc3-8aba-e6dbc1f8c911/bebff58a-a1df-4e6a-868b-ba584b296f85"]'
$ws.Cells.Item(4, 13).Value = '""'
$ws.Cells.Item(4, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "3-ZW W\u00fcrzburg  Sanderau,  3. OG.  Bad neu.", "basicObjectPricEur": 950, "basicLivingSpace": 79, "basicRooms": 3, "basicConstructionYear": 1911, "basicCity": "W\u00fcrzburg", "basicStreet": "Eichendorffstrasse 4", "basicUrl": "https://www.immowelt.de/expose/25dtq5g", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 5 <= original row 10 (url "https://www.immowelt.de/expose/2a97e5s")
$ws.Cells.Item(5, 1).Value = '"https://www.immowelt.de/expose/2a97e5s"'
$ws.Cells.Item(5, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(5, 3).Value = '"Gut gemacht !! 3 Zi, 81 qm, Balkon und Einbauk\u00fcche"'
$ws.Cells.Item(5, 4).Value = '1183.37'
$ws.Cells.Item(5, 5).Value = '"EUR"'
$ws.Cells.Item(5, 6).Value = '81.32'
$ws.Cells.Item(5, 7).Value = '3'
$ws.Cells.Item(5, 8).Value = '2019'
$ws.Cells.Item(5, 9).Value = '["Personenaufzug", "Einbauk\u00fcche", "Stellplatz"]'
$ws.Cells.Item(5, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 13", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(5, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3550", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 1183.37, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 252, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 171, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1666.37, "Unit": "EUR"}]}'
$ws.Cells.Item(5, 12).Value = '["https://ms.immowelt.org/4e4fecdb-3fe4-4b73-be6f-d02dcff6e907/a23bae17-4a7c-4155-88b0-8014e5e3021a", "https://ms.immowelt.org/324f63aa-2fff-4fc8-8ba1-ac69ea7a7e23/ed4f84b2-861b-4e99-b1c6-4e3963ae1f35", "https://ms.immowelt.org/be4e9cc8-6ad5-4ae2-a739-1b8b7178ee12/ba6c9877-b521-46d0-a222-9e1e70f8ffd4"]'
$ws.Cells.Item(5, 13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "herr", "firstName": "Stefan", "lastName": "Brogl", "phone": "+49 931 306 990 20"}'
$ws.Cells.Item(5, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Gut gemacht !! 3 Zi, 81 qm, Balkon und Einbauk\u00fcche", "basicObjectPricEur": 1183.37, "basicLivingSpace": 81.32, "basicRooms": 3, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 13", "basicUrl": "https://www.immowelt.de/expose/2a97e5s", "basicContactPhone": "+49 931 306 990 20"}'

# Row 6 <= original row 3 (url "https://www.immowelt.de/expose/2ajt75p")
$ws.Cells.Item(6, 1).Value = '"https://www.immowelt.de/expose/2ajt75p"'
$ws.Cells.Item(6, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(6, 3).Value = '"Wohnberechtigungsschein (EOF3) erforderlich!! 2-Zimmer Wohnung mit Dachterrasse!!"'
$ws.Cells.Item(6, 4).Value = '474.53'
$ws.Cells.Item(6, 5).Value = '"EUR"'
$ws.Cells.Item(6, 6).Value = '49.95'
$ws.Cells.Item(6, 7).Value = '2'
$ws.Cells.Item(6, 8).Value = '2019'
$ws.Cells.Item(6, 9).Value = '["Personenaufzug"]'
$ws.Cells.Item(6, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 15", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(6, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1420", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 474.53, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 167, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 98.9, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 740.43, "Unit": "EUR"}]}'
$ws.Cells.Item(6, 12).Value = '["https://ms.immowelt.org/059062b5-c451-4105-bfa0-dfb0db8ca29f/0f1cbcbd-e9ed-45e6-96e8-910e278b3044", "https://ms.immowelt.org/ca46e538-90d4-4fa2-8f22-2b6881ac1800/9a457fdc-44b3-4336-b20e-89b86824b294", "https://ms.immowelt.org/a1a88b02-0504-429f-b458-7c1a33126061/b8bc3487-e427-4a1b-b30b-2df51d481f90", "https://ms.immowelt.org/d507b7b4-5e3c-452e-957a-b5636948e585/d2a53e80-0b41-44e4-a681-ce8531a32c32", "https://ms.immowelt.org/ff43f21c-a328-4711-9ce2-41a9bff1bf59/83bd2e55-f9e8-44b1-bf30-e874b0a81423"]'
$ws.Cells.Item(6, 13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(6, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnberechtigungsschein (EOF3) erforderlich!! 2-Zimmer Wohnung mit Dachterrasse!!", "basicObjectPricEur": 474.53, "basicLivingSpace": 49.95, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 15", "basicUrl": "https://www.immowelt.de/expose/2ajt75p"}'

# Row 7 <= original row 4 (url "https://www.immowelt.de/expose/2tfuz4l")
$ws.Cells.Item(7, 1).Value = '"https://www.immowelt.de/expose/2tfuz4l"'
$ws.Cells.Item(7, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(7, 3).Value = '"Bestlage mit Festungsblick"'
$ws.Cells.Item(7, 4).Value = '610'
$ws.Cells.Item(7, 5).Value = '"EUR"'
$ws.Cells.Item(7, 6).Value = '65'
$ws.Cells.Item(7, 7).Value = '1'
$ws.Cells.Item(7, 8).Value = 'null'
$ws.Cells.Item(7, 9).Value = '["DUSCHE", "FENSTER", "Dachgeschoss", "frei", "offene K\u00fcche", "Stellplatz"]'
$ws.Cells.Item(7, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Sanderau", "ZipCode": "97072", "Street": "Breslauer Stra\u00dfe", "LocationId": 496022, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(7, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "Abl\u00f6se der K\u00fcche", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Sanderau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg/mietspiegel"}}, "DataTable": [{"NumberValue": 610, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete", "Comments": ["zzgl. Nebenkosten"]}]}'
$ws.Cells.Item(7, 12).Value = '["https://ms.immowelt.org/b4607fa6-96c6-40a3-8049-ed2221dafd3c/bb9b114c-9b74-4749-b16f-f5402ce89d13", "https://ms.immowelt.org/0678cd9c-641a-4bcc-a0ac-02f46fc8e1ea/b0aa6847-a6e4-474c-b300-23bb59c7628b", "https://ms.immowelt.org/d524e2d3-1330-4559-987d-6495b2a4c71d/3adebae2-58ee-44b1-a898-86e63f8eff12", "https://ms.immowelt.org/61eb287f-36c7-4a51-a4db-f88e4049ad30/29721137-7a70-4a6d-9abc-340c019663bf"]'
$ws.Cells.Item(7, 13).Value = '""'
$ws.Cells.Item(7, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Bestlage mit Festungsblick", "basicObjectPricEur": 610, "basicLivingSpace": 65, "basicRooms": 1, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicStreet": "Breslauer Stra\u00dfe", "basicUrl": "https://www.immowelt.de/expose/2tfuz4l", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 8 <= original row 12 (url "https://www.immowelt.de/expose/2aqpt5d")
$ws.Cells.Item(8, 1).Value = '"https://www.immowelt.de/expose/2aqpt5d"'
$ws.Cells.Item(8, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(8, 3).Value = '"Die Einbauk\u00fcche ist schon drin - Sie m\u00fcssen nur noch einziehen! Gehobene 2-Z. Wohnung im Steinbachtal "'
$ws.Cells.Item(8, 4).Value = '1165'
$ws.Cells.Item(8, 5).Value = '"EUR"'
$ws.Cells.Item(8, 6).Value = '77'
$ws.Cells.Item(8, 7).Value = '2'
$ws.Cells.Item(8, 8).Value = 'null'
$ws.Cells.Item(8, 9).Value = '["renoviert", "GAS", "Zentralheizung", "gartennutzung", "Einbauk\u00fcche", "Terrasse"]'
$ws.Cells.Item(8, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Steinbachtal", "ZipCode": "97082", "LocationId": 496023, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(8, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3495", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Steinbachtal)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-steinbachtal/mietspiegel"}}, "DataTable": [{"NumberValue": 1165, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 280, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1445, "Unit": "EUR"}, {"NumberValue": 75, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "1 Stellplatz"}]}'
$ws.Cells.Item(8, 12).Value = '["https://ms.immowelt.org/d789ca7c-5451-48cf-a649-7af7d4c9078a/0af4290a-753d-4df9-bd0e-3989d2db4215", "https://ms.immowelt.org/2908c42c-0f96-410b-b619-dba38e6a42a1/75b141df-d3a1-4754-b265-67306b3a36e1", "https://ms.immowelt.org/fbb8e8e1-0cc5-4a95-8902-c6f6af96acc8/548e065f-569a-4838-9887-a19794c56f49", "https://ms.immowelt.org/4365c939-93fc-4ec1-b3e8-77d586d0e200/5112a29b-5171-4dc0-b36a-b9b4c4ac411f", "https://ms.immowelt.org/db7ac2e7-2969-42be-920b-36be7e6aab0d/0cf7253b-4260-4efa-af51-0074b8782be6", "https://ms.immowelt.org/d1a4bcbb-b277-4731-8486-0a09fb1b8996/475444d0-b794-44bc-ba59-ce3f2f51b9f8", "https://ms.immowelt.org/748a1662-4a4a-4f6d-bb72-baf731de8adc/700fbb9a-6ef9-48a1-8e73-4c46cba5c258", "https://ms.immowelt.org/ed9af812-426e-4e49-8e89-2ede0c407269/5a107fdf-e72a-4f22-9de2-52351751a755", "https://ms.immowelt.org/643489b1-7bf2-46a3-882c-7d8b8b0957ca/bce02d3c-331a-4110-8e99-d872140474b1", "https://ms.immowelt.org/d8490c94-8bc2-4f38-affe-7986af339d1c/5d00296c-3134-4d05-a29a-760985379123", "https://ms.immowelt.org/862887dd-fca5-4504-a254-d22277945d21/1fdaef5f-474b-4b44-84df-c23fb4960205", "https://ms.immowelt.org/235cb2e7-3834-406d-a4ac-0a1db7702757/b3d2810f-2282-453d-99bf-f70073155c71", "https://ms.immowelt.org/3ae43d7c-7b0d-435e-84bb-7a6ea4aafc03/ae2064f4-8fe0-463b-958f-a7c34e88466a"]'
$ws.Cells.Item(8, 13).Value = '""'
$ws.Cells.Item(8, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Die Einbauk\u00fcche ist schon drin - Sie m\u00fcssen nur noch einziehen! Gehobene 2-Z. Wohnung im Steinbachtal ", "basicObjectPricEur": 1165, "basicLivingSpace": 77, "basicRooms": 2, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2aqpt5d", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 9 <= original row 5 (url "https://www.immowelt.de/expose/2aaj75n")
$ws.Cells.Item(9, 1).Value = '"https://www.immowelt.de/expose/2aaj75n"'
$ws.Cells.Item(9, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(9, 3).Value = '"Wohnung kann nur mit Wohnberechtigungsschein angemietet werden! 2 Zimmer auf 56 m\u00b2"'
$ws.Cells.Item(9, 4).Value = '532.19'
$ws.Cells.Item(9, 5).Value = '"EUR"'
$ws.Cells.Item(9, 6).Value = '56.02'
$ws.Cells.Item(9, 7).Value = '2'
$ws.Cells.Item(9, 8).Value = '2019'
$ws.Cells.Item(9, 9).Value = '["Personenaufzug"]'
$ws.Cells.Item(9, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 17", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(9, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1590", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 532.19, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 179, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 110.92, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 882.11, "Unit": "EUR"}]}'
$ws.Cells.Item(9, 12).Value = '["https://ms.immowelt.org/814cec20-d48a-4483-ae7f-5b779e29c323/404d9d3f-baa2-473d-af13-096d09772807", "https://ms.immowelt.org/60ba3310-b6df-446e-b860-6deb11ef0a09/2558c572-051e-464d-9aad-205683e0f215", "https://ms.immowelt.org/e10fa924-843a-4cd9-897d-8a60a9cd8473/3fd9082b-7a00-4210-99b9-f828a6a2e085", "https://ms.immowelt.org/20e41405-47ae-47ac-97d8-0e7340de1982/b112096a-8211-4248-a5e6-0904121ddac0", "https://ms.immowelt.org/2300c3d0-702d-4867-94fc-810ec8ca8971/f6adb669-5e8c-4182-ad4f-a5b6f477b025", "https://ms.immowelt.org/bdd38a88-61e2-4f96-bac0-895a2e862eed/9b1e1d00-b895-496c-bf47-2912364590ab", "https://ms.immowelt.org/a047e5df-c99a-4420-a028-e84ef22f3786/0dd57565-cb4f-4092-ab20-2a61b660568a"]'
$ws.Cells.Item(9, 13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(9, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnung kann nur mit Wohnberechtigungsschein angemietet werden! 2 Zimmer auf 56 m\u00b2", "basicObjectPricEur": 532.19, "basicLivingSpace": 56.02, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 17", "basicUrl": "https://www.immowelt.de/expose/2aaj75n"}'

# Row 10 <= original row 11 (url "https://www.immowelt.de/expose/2afyn5t")
$ws.Cells.Item(10, 1).Value = '"https://www.immowelt.de/expose/2afyn5t"'
$ws.Cells.Item(10, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(10, 3).Value = '"2,5 Zimmer Wohnung beim Japanischen Garten in W\u00fcrzburg"'
$ws.Cells.Item(10, 4).Value = '750'
$ws.Cells.Item(10, 5).Value = '"EUR"'
$ws.Cells.Item(10, 6).Value = '65'
$ws.Cells.Item(10, 7).Value = '2.5'
$ws.Cells.Item(10, 8).Value = '1960'
$ws.Cells.Item(10, 9).Value = '["Altbau (bis 1945)", "WANNE", "Kelleranteil", "renoviert", "GAS", "Zentralheizung", "Einbauk\u00fcche"]'
$ws.Cells.Item(10, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Zellerau", "ZipCode": "97082", "LocationId": 496026, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(10, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "2000", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Zellerau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-zellerau/mietspiegel"}}, "DataTable": [{"NumberValue": 750, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 200, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 950, "Unit": "EUR"}]}'
$ws.Cells.Item(10, 12).Value = '["https://ms.immowelt.org/7392265b-7589-4e96-8719-b289bb8ace01/19ed71d6-701a-4bbe-929c-a35f0d2ec94b", "https://ms.immowelt.org/15baddee-aaf7-4f0a-b016-7f47a07c4e37/48316143-5143-447a-8c57-f0b7e929d6fe", "https://ms.immowelt.org/027c0258-c7df-4398-b7bc-fce818803786/c5879894-c70e-4d36-87c5-2f3d1391a860", "https://ms.immowelt.org/512f5ff2-4810-463d-be09-2b1645c65d7d/c20634e3-424e-4fa7-adc1-7bb96e05a731", "https://ms.immowelt.org/0bc625b6-c3c3-4417-81a8-a0c000046701/97c3d523-e7c6-4b79-af4f-d9928909527c", "https://ms.immowelt.org/76a5fc74-de02-4ae9-aec8-48e4168cea64/423affd6-d43f-46c2-9187-ca30a02614d6", "https://ms.immowelt.org/39f67b57-3390-44ea-ab71-3cbb1b166a73/03cd65f9-4b01-46f7-832d-c4ec60e69999", "https://ms.immowelt.org/4aa9aeb4-7027-43a0-85f3-887e51080c31/b7c1b3ae-576c-44d3-a98d-260591a7c97b"]'
$ws.Cells.Item(10, 13).Value = '""'
$ws.Cells.Item(10, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "2,5 Zimmer Wohnung beim Japanischen Garten in W\u00fcrzburg", "basicObjectPricEur": 750, "basicLivingSpace": 65, "basicRooms": 2.5, "basicConstructionYear": 1960, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2afyn5t", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 11 <= original row 9 (url "https://www.immowelt.de/expose/2a7j75n")
$ws.Cells.Item(11, 1).Value = '"https://www.immowelt.de/expose/2a7j75n"'
$ws.Cells.Item(11, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(11, 3).Value = '"Wohnen im herrlichen HUBland: 2Zi, 70m\u00b2, EBK und Balkon"'
$ws.Cells.Item(11, 4).Value = '1181.21'
$ws.Cells.Item(11, 5).Value = '"EUR"'
$ws.Cells.Item(11, 6).Value = '70.87'
$ws.Cells.Item(11, 7).Value = '2'
$ws.Cells.Item(11, 8).Value = '2019'
$ws.Cells.Item(11, 9).Value = '["Personenaufzug", "Einbauk\u00fcche", "Stellplatz"]'
$ws.Cells.Item(11, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 13", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(11, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3540", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 1181.21, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 219, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 169, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1629.21, "Unit": "EUR"}]}'
$ws.Cells.Item(11, 12).Value = '["https://ms.immowelt.org/0d9ec1ef-21cc-484a-88ae-193c7267a3bb/1e78105c-6490-4468-a6c3-936b8b388b80", "https://ms.immowelt.org/5c55f670-770f-44e9-b26b-7c9164f4a76a/9f56e225-8915-45d3-9af7-0e36ad2d804e", "https://ms.immowelt.org/12a9fd2a-da00-49a5-9999-efb246a9d727/efc35e4b-a860-47d1-823d-1282d60e0bc6", "https://ms.immowelt.org/bb9a6b97-0d40-46d1-9b24-011a5dcc7508/40526838-218d-4a4c-9422-8bd98023d9f6", "https://ms.immowelt.org/1a701923-58c4-4074-92b3-3e6bc100d51c/2f4849ee-2b95-432d-8ed4-b41b0211a117", "https://ms.immowelt.org/4d6aee50-6e8f-423f-b67a-bd01a011592d/4d58cead-e166-4a77-8c0f-2ee869a9e824"]'
$ws.Cells.Item(11, 13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "herr", "firstName": "Stefan", "lastName": "Brogl", "phone": "+49 931 306 990 20"}'
$ws.Cells.Item(11, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnen im herrlichen HUBland: 2Zi, 70m\u00b2, EBK und Balkon", "basicObjectPricEur": 1181.21, "basicLivingSpace": 70.87, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 13", "basicUrl": "https://www.immowelt.de/expose/2a7j75n", "basicContactPhone": "+49 931 306 990 20"}'

# Row 12 <= original row 13 (url "https://www.immowelt.de/expose/2at6f5v")
$ws.Cells.Item(12, 1).Value = '"https://www.immowelt.de/expose/2at6f5v"'
$ws.Cells.Item(12, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(12, 3).Value = '"Wohnung vermieten in W\u00fcrzburg"'
$ws.Cells.Item(12, 4).Value = '420'
$ws.Cells.Item(12, 5).Value = '"EUR"'
$ws.Cells.Item(12, 6).Value = '""'
$ws.Cells.Item(12, 7).Value = '""'
$ws.Cells.Item(12, 8).Value = 'null'
$ws.Cells.Item(12, 9).Value = '["GAS"]'
$ws.Cells.Item(12, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Gromb\u00fchl", "ZipCode": "97080", "LocationId": 496011, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(12, 11).Value = '{"AdditionalInformation": {"MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Gromb\u00fchl)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-grombuehl/mietspiegel"}}, "DataTable": [{"NumberValue": 420, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete", "Comments": ["zzgl. Nebenkosten"]}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 420, "Unit": "EUR"}]}'
$ws.Cells.Item(12, 12).Value = '[]'
$ws.Cells.Item(12, 13).Value = '""'
$ws.Cells.Item(12, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnung vermieten in W\u00fcrzburg", "basicObjectPricEur": 420, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2at6f5v", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 13 <= original row 15 (url "https://www.immowelt.de/expose/2asr65p")
$ws.Cells.Item(13, 1).Value = '"https://www.immowelt.de/expose/2asr65p"'
$ws.Cells.Item(13, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(13, 3).Value = '"Familienfreundliche 4 - Zimmer Wohnung in W\u00fcrzburg - Heidingsfeld "'
$ws.Cells.Item(13, 4).Value = '1100'
$ws.Cells.Item(13, 5).Value = '"EUR"'
$ws.Cells.Item(13, 6).Value = '103'
$ws.Cells.Item(13, 7).Value = '4'
$ws.Cells.Item(13, 8).Value = '1984'
$ws.Cells.Item(13, 9).Value = '["WANNE", "gaestewc", "FENSTER", "Kelleranteil", "renoviert", "GAS", "Etagenheizung", "gartennutzung", "Einbauk\u00fcche", "wg_geeignet", "Balkon"]'
$ws.Cells.Item(13, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Heidingsfeld", "ZipCode": "97084", "Street": "Unterer Kirchbergweg 37", "LocationId": 496013, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(13, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "2200,00", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Heidingsfeld)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-heidingsfeld/mietspiegel"}}, "DataTable": [{"NumberValue": 1100, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 300, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1400, "Unit": "EUR"}, {"NumberValue": 70, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "1 Stellplatz"}]}'
$ws.Cells.Item(13, 12).Value = '["https://ms.immowelt.org/bbd4dcbf-c64e-4c60-9789-afc66db53a3f/7fc04ce4-7f84-4ce5-8111-31fcc917306c", "https://ms.immowelt.org/315ff83e-d234-453a-ad52-c4b99bef7055/a979f182-320d-49d1-8002-bcde353d36ee", "https://ms.immowelt.org/6137c704-6a29-47ae-a5a3-88abd0a2fde0/d535204c-f43a-4229-8152-507fc644a211", "https://ms.immowelt.org/1981c2ad-23cb-41f4-825a-77d67f5fbbe7/e968ecc9-ea46-4188-a2c4-968720f069c0", "https://ms.immowelt.org/a4965ee1-de2b-4d90-8333-b793f90a7406/867d9997-6f27-48da-ade0-9096bf31cf3e", "https://ms.immowelt.org/aed0a388-445d-4f8f-9afd-4a146b35d962/3bc7e507-92a3-446a-a3f1-5799c1df7338", "https://ms.immowelt.org/f9d8ba2e-fb1a-4cdd-98ea-ba6ba4973c25/fde9d74e-bd14-4413-a148-81aad236c43f", "https://ms.immowelt.org/e0fbe89c-ecbb-41e9-9ecd-61c2762b9189/3784b2d8-f66d-4f7f-a893-cb9c7b3e982c", "https://ms.immowelt.org/16af38b3-f1c4-4d16-ba17-2685efc0dbe6/b5361690-0608-4e15-8167-e5ab29ac41e7", "https://ms.immowelt.org/fb842589-0c99-45bb-b567-5d04f811adf5/729e7fb7-c4eb-47db-936a-8a406b3efda8", "https://ms.immowelt.org/2f97e2d2-4fcd-40db-bee9-5d4161812bb0/a57df980-0cf8-4352-8cf0-5e867376c343", "https://ms.immowelt.org/79cf67c5-73a5-4097-8b85-e20df2ab786e/363fc58b-d0d5-4f6e-afda-b1d880cf8747", "https://ms.immowelt.org/d9ae7b88-e98c-4d7f-8096-8ede0076f326/e28b854c-9b59-46d6-84a5-28d2595f9ff2", "https://ms.immowelt.org/dffe6319-03a2-4071-98a5-e0f1e3954460/de63d734-7d33-492b-acbf-be5ca5b61e49", "https://ms.immowelt.org/8af5caf8-2627-4044-83b1-a3824003abcd/a9a95b7e-658f-4e56-989c-5a564997561b", "https://ms.immowelt.org/de373ea2-fa91-495b-aa05-ffa24233b025/ff4e0de7-6c79-411b-bec2-3f5b973c460c", "https://ms.immowelt.org/1db3601f-db1c-4e6f-8b1d-e2f17c88804b/e84274ae-8668-49f0-8327-447ff8511752", "https://ms.immowelt.org/6bc525ad-0769-42c3-aaff-a169a251a660/789b8d1b-633d-4989-a0e8-50949ea51d51"]'
$ws.Cells.Item(13, 13).Value = '{"salutation": "Herr", "firstName": "Arnold", "lastName": "Hammer", "mobile": "017683299994", "phone": "0931/60565"}'
$ws.Cells.Item(13, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Familienfreundliche 4 - Zimmer Wohnung in W\u00fcrzburg - Heidingsfeld ", "basicObjectPricEur": 1100, "basicLivingSpace": 103, "basicRooms": 4, "basicConstructionYear": 1984, "basicCity": "W\u00fcrzburg", "basicStreet": "Unterer Kirchbergweg 37", "basicUrl": "https://www.immowelt.de/expose/2asr65p", "basicContactPhone": "0931/60565", "basicContactMobile": "017683299994"}'

# Row 14 <= original row 7 (url "https://www.immowelt.de/expose/2a8j75n")
$ws.Cells.Item(14, 1).Value = '"https://www.immowelt.de/expose/2a8j75n"'
$ws.Cells.Item(14, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(14, 3).Value = '"Wohnberechtigungsschein erforderlich - 2 Zi, 58m\u00b2, Balkon"'
$ws.Cells.Item(14, 4).Value = '554.99'
$ws.Cells.Item(14, 5).Value = '"EUR"'
$ws.Cells.Item(14, 6).Value = '58.42'
$ws.Cells.Item(14, 7).Value = '2'
$ws.Cells.Item(14, 8).Value = '2019'
$ws.Cells.Item(14, 9).Value = '["Personenaufzug", "Stellplatz"]'
$ws.Cells.Item(14, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Frauenland", "ZipCode": "97074", "Street": "Athanasius-Kircher-Stra\u00dfe 15", "LocationId": 496008, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(14, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1660", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Frauenland)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-frauenland/mietspiegel"}}, "DataTable": [{"NumberValue": 554.99, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 210, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 208, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1032.99, "Unit": "EUR"}]}'
$ws.Cells.Item(14, 12).Value = '["https://ms.immowelt.org/2c287869-b002-4a19-baea-0433bdc7a327/e15bc071-fe33-4461-8b34-9efe6f7f6e23", "https://ms.immowelt.org/5b165c7f-f365-48c1-841e-a8b1fd2bb6a1/23a21cb8-5f33-4e46-b102-f4960458c3e9", "https://ms.immowelt.org/4b517be9-8d93-492f-ac56-516e6bad3d95/1238a7b6-1272-4e53-9ead-1eafab0548c7", "https://ms.immowelt.org/acdfaf67-7ef9-4ee5-a2e3-416020180ea6/b120a8ca-7b3e-47f0-99b1-006c6eb00277", "https://ms.immowelt.org/540c5813-745e-419c-9217-36ba5261e842/259d75f3-7aca-466d-b4ed-f4b5a08ff770", "https://ms.immowelt.org/9ac46b17-c8c4-446e-8e96-de026ce267b4/f6098bb6-68df-4742-8eda-56e4660c3221", "https://ms.immowelt.org/09537047-f692-4086-b911-c08acb726d74/595eae7d-cd2e-4c30-b4d5-6bf58c421952"]'
$ws.Cells.Item(14, 13).Value = '{"companyName": "BUWOG Immobilien Treuhand GmbH\u00ad", "address": {"city": "Kiel", "zipCode": "24103", "street": "Fabrikstra\u00dfe 7"}, "salutation": "Herr", "firstName": "Stefan", "lastName": "Brogl"}'
$ws.Cells.Item(14, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Wohnberechtigungsschein erforderlich - 2 Zi, 58m\u00b2, Balkon", "basicObjectPricEur": 554.99, "basicLivingSpace": 58.42, "basicRooms": 2, "basicConstructionYear": 2019, "basicCity": "W\u00fcrzburg", "basicStreet": "Athanasius-Kircher-Stra\u00dfe 15", "basicUrl": "https://www.immowelt.de/expose/2a8j75n"}'

# Row 15 <= original row 16 (url "https://www.immowelt.de/expose/2awv45u")
$ws.Cells.Item(15, 1).Value = '"https://www.immowelt.de/expose/2awv45u"'
$ws.Cells.Item(15, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(15, 3).Value = '"Luxuri\u00f6se Penthouse-Wohnung Erstbezug! 3-Zimmer-Neubau in W\u00fcrzburg-Heidingfeld"'
$ws.Cells.Item(15, 4).Value = '1720'
$ws.Cells.Item(15, 5).Value = '"EUR"'
$ws.Cells.Item(15, 6).Value = '123.36'
$ws.Cells.Item(15, 7).Value = '3'
$ws.Cells.Item(15, 8).Value = '2022'
$ws.Cells.Item(15, 9).Value = '["Neubau", "barriefrei", "WANNE", "gaestewc", "FENSTER", "Kelleranteil", "Personenaufzug", "FERNE", "GAS", "Fu\u00dfbodenheizung", "Zentralheizung", "rollstuhlgerecht", "PARKETT", "frei", "Einbauk\u00fcche", "offene K\u00fcche", "Speisekammer", "Stellplatz", "Tiefgarage", "REINIGUNG", "Balkon", "Terrasse"]'
$ws.Cells.Item(15, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Heidingsfeld", "ZipCode": "97084", "LocationId": 496013, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(15, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "5160", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Heidingsfeld)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-heidingsfeld/mietspiegel"}}, "DataTable": [{"NumberValue": 1720, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 320, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 2040, "Unit": "EUR"}, {"NumberValue": 70, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "2 Stellpl\u00e4tze, je"}]}'
$ws.Cells.Item(15, 12).Value = '["https://ms.immowelt.org/70edfd50-f060-4619-8d8e-0a42b4111160/271e31b7-8484-45a8-8ac9-904dabf13c7b", "https://ms.immowelt.org/efa2a84b-bb24-4421-ae15-0649647fc5f0/8490cb37-3682-446f-8b11-e9e1da37cba9", "https://ms.immowelt.org/1683eefc-9ca0-4cbc-992f-8e3a9b044b29/7df67d66-b166-4c09-b89c-52d2409839a6", "https://ms.immowelt.org/80c22a83-a553-4a11-b6d3-845ab5aee5fa/a5ed848d-87ec-43a7-91b1-3a82889e2f11", "https://ms.immowelt.org/912ad1a7-5c04-49c1-bfdb-52dfd32f4078/087254f6-99be-4f33-ad19-a0d31350b6df", "https://ms.immowelt.org/c04aac07-656d-484f-97c8-f9183167c740/27c413f7-27b3-4874-9d22-7b3ce7fb1772", "https://ms.immowelt.org/6355ac0c-0919-41d1-9939-c7327ad98472/b5a33c7f-e707-4bd6-839a-ea939f810d7c", "https://ms.immowelt.org/e86eec55-1463-4303-8c06-6e1cd90ac14e/01d5e414-8279-49e7-ac42-160f28ed6a9b", "https://ms.immowelt.org/f6c4b53e-148d-44ad-b8fb-50e2e7c095c6/9403fef4-e110-4d99-88ec-08c732282eda", "https://ms.immowelt.org/bb63f6b4-e5b8-4f84-b1e9-c19308a3d41f/20e645fd-f306-449c-82ae-d3591730cca2", "https://ms.immowelt.org/6eff25e1-ca40-4fdc-b6f2-adcf3b4d0a45/63b7a19b-f618-4933-adb0-f0cd5e1263a5", "https://ms.immowelt.org/2811d6e1-5985-40bf-a07a-594fdaf52254/4a7a045f-966f-4645-86ce-900f6e22d446", "https://ms.immowelt.org/aefa4aed-f85d-4406-ae1b-9c233616a69e/fd830bae-b6de-4525-8e19-1c72c8249092", "https://ms.immowelt.org/c0855678-63ef-4feb-a503-a8d0adf77171/3851bcbc-db22-481c-b205-d0b9fb09017f", "https://ms.immowelt.org/df7e819e-aa27-4bbf-b743-a9e2712de66b/0737cd62-424e-4825-afba-0a294eb938ec", "https://ms.immowelt.org/f47e219f-f548-45dd-a07d-bb94c2b60e97/cfc0cf18-e610-4657-b7ed-3ca18984267a"]'
$ws.Cells.Item(15, 13).Value = '{"companyName": "Gute Bude Immobilien", "address": {"city": "W\u00fcrzburg", "zipCode": "97082", "street": "Zeller Str. 3c"}, "salutation": "Frau", "firstName": "Nelly", "lastName": "Gronau"}'
$ws.Cells.Item(15, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Luxuri\u00f6se Penthouse-Wohnung Erstbezug! 3-Zimmer-Neubau in W\u00fcrzburg-Heidingfeld", "basicObjectPricEur": 1720, "basicLivingSpace": 123.36, "basicRooms": 3, "basicConstructionYear": 2022, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2awv45u"}'

# Row 16 <= original row 20 (url "https://www.immowelt.de/expose/2upu84s")
$ws.Cells.Item(16, 1).Value = '"https://www.immowelt.de/expose/2upu84s"'
$ws.Cells.Item(16, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(16, 3).Value = '"Gro\u00dfe 3-Zimmer-Wohnung mit eigenem Garten und Balkon in K\u00fcrnach"'
$ws.Cells.Item(16, 4).Value = '900'
$ws.Cells.Item(16, 5).Value = '"EUR"'
$ws.Cells.Item(16, 6).Value = '90'
$ws.Cells.Item(16, 7).Value = '3'
$ws.Cells.Item(16, 8).Value = 'null'
$ws.Cells.Item(16, 9).Value = '["WANNE", "FENSTER", "Kelleranteil", "ELEKTRO", "Garten", "gartennutzung", "Balkon", "Terrasse"]'
$ws.Cells.Item(16, 10).Value = '{"Country": "Deutschland", "City": "K\u00fcrnach", "ZipCode": "97273", "LocationId": 12666, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(16, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "2 Kaltmieten ", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in K\u00fcrnach", "Link": "https://www.immowelt.de/immobilienpreise/kuernach/mietspiegel"}}, "DataTable": [{"NumberValue": 900, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 80, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 50, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "1 Stellplatz"}]}'
$ws.Cells.Item(16, 12).Value = '["https://ms.immowelt.org/11008e51-5945-4ab3-a0eb-f2a703785b20/53fa6be7-ade6-4395-81ca-d66b9bbb1ba9", "https://ms.immowelt.org/64bc93db-da55-44ca-a2c8-c14f8ce86353/14ba8183-9d5b-4c81-a957-284d595a032f", "https://ms.immowelt.org/87e6208f-ee9e-4d96-a62a-6561d184b0e3/572c15dc-ed2e-4296-bc44-96a28e7cf135", "https://ms.immowelt.org/ba824adc-5910-4328-b432-ac6c69a928b3/dc0e3963-adad-4942-ad9a-dc0fc1b71796", "https://ms.immowelt.org/cee10f88-e7a3-499e-be76-626c1d42e77d/c1805001-8c80-4257-977c-80d06c2cbbc7", "https://ms.immowelt.org/aba27baa-44fb-49d2-833b-ffe9285cf719/2f9875be-ad3f-435e-bf3a-d7d63a5e8a17", "https://ms.immowelt.org/b63c1740-9745-464b-8364-2df590a1d676/f71392e8-2cac-4fad-aeaf-0379e7187348", "https://ms.immowelt.org/0b219ac9-8c45-4c18-9228-d206e9ff703e/4bce456f-aecf-48fa-a4a2-b0b332e7badc", "https://ms.immowelt.org/711ca5fb-12d7-42dc-a1f6-32893287c0eb/4002136b-ce99-43ba-8406-dee393492ff0", "https://ms.immowelt.org/17453203-a015-4bcc-816b-12fe9fe761ee/3c33ba0b-8a93-4d03-b3c7-0267e690bc76", "https://ms.immowelt.org/f09d3348-84eb-4a67-941e-68eee061b36d/75a81b43-153a-44e2-8336-9f323ccdad5a", "https://ms.immowelt.org/79e28d05-fb4a-4b73-ad94-cd78318ef65b/0441105e-ad78-4d19-86b1-dab1c3c87ee1", "https://ms.immowelt.org/07c2164c-067e-490d-b769-760ef64271e2/2e82cebf-00da-4337-83e8-6186da964d0d", "https://ms.immowelt.org/de5a7379-e440-43c4-8670-4bdd1e8bdca0/e80fa6cc-bbb1-48ab-b091-2b052e63948b", "https://ms.immowelt.org/9cc41d7e-f16c-4d38-ae23-e1f15da29cfc/e1740437-8eb5-4984-bba2-1c1036d06b98"]'
$ws.Cells.Item(16, 13).Value = '""'
$ws.Cells.Item(16, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Gro\u00dfe 3-Zimmer-Wohnung mit eigenem Garten und Balkon in K\u00fcrnach", "basicObjectPricEur": 900, "basicLivingSpace": 90, "basicRooms": 3, "basicConstructionYear": null, "basicCity": "K\u00fcrnach", "basicUrl": "https://www.immowelt.de/expose/2upu84s", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 17 <= original row 14 (url "https://www.immowelt.de/expose/2agbf5s")
$ws.Cells.Item(17, 1).Value = '"https://www.immowelt.de/expose/2agbf5s"'
$ws.Cells.Item(17, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(17, 3).Value = '"Traumwohnung in Top Lage mit toller Weitsicht!"'
$ws.Cells.Item(17, 4).Value = '2880'
$ws.Cells.Item(17, 5).Value = '"EUR"'
$ws.Cells.Item(17, 6).Value = '195'
$ws.Cells.Item(17, 7).Value = '6'
$ws.Cells.Item(17, 8).Value = '2023'
$ws.Cells.Item(17, 9).Value = '["Neubau", "barriefrei", "WANNE", "gaestewc", "FENSTER", "Kelleranteil", "LUFTWP", "Garten", "gartennutzung", "rollstuhlgerecht", "Balkon", "Terrasse"]'
$ws.Cells.Item(17, 10).Value = '{"Country": "Deutschland", "City": "Gerbrunn", "ZipCode": "97218", "Street": "Helene Wessel Weg 5", "LocationId": 12047, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(17, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "3 Kaltmieten", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in Gerbrunn", "Link": "https://www.immowelt.de/immobilienpreise/gerbrunn/mietspiegel"}}, "DataTable": [{"NumberValue": 2880, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 400, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 3280, "Unit": "EUR"}]}'
$ws.Cells.Item(17, 12).Value = '["https://ms.immowelt.org/17e5aca9-793d-40b3-aab5-a5dcea349a87/0d9236ef-f790-4bac-82dc-9a229cfe4527", "https://ms.immowelt.org/cf96b674-5a06-4156-b363-d3e23b5dfe92/c12ffbff-3538-4b01-a604-ee0869a952b6", "https://ms.immowelt.org/2755511b-7b13-42f8-bad0-3d769da9305b/b3d17421-aa79-4be1-83cb-ae6fe7955297", "https://ms.immowelt.org/17109eb9-b8cf-435e-b0b9-8377b8a7c9f4/986295b9-a109-4924-b6b8-fc375860cd30", "https://ms.immowelt.org/d4384af4-fbdc-4807-9ab2-e4c02ad5e5b9/61dd1708-3529-4a40-9890-401b9917edfc", "https://ms.immowelt.org/efa4e805-be82-4e38-998b-46485660f4e5/fbb8ea8f-2b17-45b1-8691-eb2fe09dd0a7"]'
$ws.Cells.Item(17, 13).Value = '""'
$ws.Cells.Item(17, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Traumwohnung in Top Lage mit toller Weitsicht!", "basicObjectPricEur": 2880, "basicLivingSpace": 195, "basicRooms": 6, "basicConstructionYear": 2023, "basicCity": "Gerbrunn", "basicStreet": "Helene Wessel Weg 5", "basicUrl": "https://www.immowelt.de/expose/2agbf5s", "basicContactPhone": "", "basicContactMobile": ""}'

# Row 18 <= original row 19 (url "https://www.immowelt.de/expose/2a2yn5v")
$ws.Cells.Item(18, 1).Value = '"https://www.immowelt.de/expose/2a2yn5v"'
$ws.Cells.Item(18, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(18, 3).Value = '"Optimal f\u00fcr Studenten-WG : 2-Zimmerwohnung in zentraler Lage im Frauenland"'
$ws.Cells.Item(18, 4).Value = '580'
$ws.Cells.Item(18, 5).Value = '"EUR"'
$ws.Cells.Item(18, 6).Value = '44'
$ws.Cells.Item(18, 7).Value = '2'
$ws.Cells.Item(18, 8).Value = 'null'
$ws.Cells.Item(18, 9).Value = '["WANNE"]'
$ws.Cells.Item(18, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Altstadt", "ZipCode": "97072", "LocationId": 496004, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(18, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "1.080,00 ", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Altstadt)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-altstadt/mietspiegel"}}, "DataTable": [{"NumberValue": 580, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 120, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Warmmiete enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}]}'
$ws.Cells.Item(18, 12).Value = '["https://ms.immowelt.org/67968400-d13b-449c-8f0f-1f8dc3645f4c/4d6bfee6-a4bb-438e-a524-96840bd51ae1", "https://ms.immowelt.org/c3319402-5904-4d76-95d2-db5b70bee409/4e68c2f5-c18e-4b63-bc8f-0ee1c1309ebf", "https://ms.immowelt.org/2dd83638-fb99-4566-ba6f-cf7b62d7c4bb/6ceb0ef0-07a4-468f-8f2a-f2e6741a577f", "https://ms.immowelt.org/1609ec66-fe8f-4f74-b4d5-02421f7513b0/a1757c96-7393-44bc-bf2e-6c9466ce4f45", "https://ms.immowelt.org/35e2e29e-1627-4973-99d4-4c3b1facdc4e/b477018c-9c03-4c1b-861c-11a14bcff4e8", "https://ms.immowelt.org/48e33aa4-fbdd-48d1-ac41-d25502baa7d5/14df5c1a-b256-4cd0-a73b-72d4dd9a2961"]'
$ws.Cells.Item(18, 13).Value = '{"companyName": "Viefhaus Immobilien Inh. Sebastian Viefhaus", "address": {"city": "W\u00fcrzburg", "zipCode": "97074", "street": "Fichtestra\u00dfe 10"}, "salutation": "Herr", "firstName": "Thomas", "lastName": "Meister", "phone": "0931-88065061"}'
$ws.Cells.Item(18, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Optimal f\u00fcr Studenten-WG : 2-Zimmerwohnung in zentraler Lage im Frauenland", "basicObjectPricEur": 580, "basicLivingSpace": 44, "basicRooms": 2, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2a2yn5v", "basicContactPhone": "0931-88065061"}'

# Row 19 <= original row 17 (url "https://www.immowelt.de/expose/2uewf4g")
$ws.Cells.Item(19, 1).Value = '"https://www.immowelt.de/expose/2uewf4g"'
$ws.Cells.Item(19, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(19, 3).Value = '"Boutique Wohnung m\u00f6bliert mieten"'
$ws.Cells.Item(19, 4).Value = '1100'
$ws.Cells.Item(19, 5).Value = '"EUR"'
$ws.Cells.Item(19, 6).Value = '55'
$ws.Cells.Item(19, 7).Value = '""'
$ws.Cells.Item(19, 8).Value = 'null'
$ws.Cells.Item(19, 9).Value = '[]'
$ws.Cells.Item(19, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Altstadt", "ZipCode": "97072", "LocationId": 496004, "PublishStreet": false, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(19, 11).Value = '{"AdditionalInformation": {"MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Altstadt)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg-altstadt/mietspiegel"}}, "DataTable": [{"NumberValue": 1100, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 200, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Warmmiete enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}]}'
$ws.Cells.Item(19, 12).Value = '["https://ms.immowelt.org/a13c7ad0-2b70-4b4d-a3ea-a34885d87b08/ec7d0988-cc69-4927-b78a-e0331046ca91", "https://ms.immowelt.org/fde8fa3b-e4a0-4ded-a7b8-8420d1a01a9f/d1572cf7-7d4c-4a02-98e6-2d6e2eaf2b62", "https://ms.immowelt.org/44adf4be-fce6-455d-9181-c27144d0b47e/2f83da82-fc2a-4b88-be4c-03e254e651ef", "https://ms.immowelt.org/ecad88ed-f5a2-43ee-9dac-2470c0ec1f30/b45f2c26-817a-4b2e-8c51-b6b4e1ffdb1d", "https://ms.immowelt.org/38e4f785-7a20-46ce-bc01-3909054d335a/ad336fba-e265-4254-b66b-537222384071", "https://ms.immowelt.org/1e7916ab-6c4d-4283-a5fa-b6ccb79f8f6e/58a84848-235f-402d-a541-f5f6b47652ab", "https://ms.immowelt.org/1160b145-63b1-49fd-9bc5-54104f7bf183/20f3af07-2776-4471-981e-ff09df9d66a7", "https://ms.immowelt.org/1498b260-43eb-4602-8d7d-ad06837fbc82/adf6d355-f74c-45b5-bfa8-b60c2f35c3b9", "https://ms.immowelt.org/3606502e-37fa-4430-9bd8-fe84c0031264/e4a61da6-7999-4ec3-b472-e1f89796a906", "https://ms.immowelt.org/d2aadef4-121e-4c0a-8d0a-767393cc0b74/6210ad93-fbc0-4207-8cf8-9c21f2312742", "https://ms.immowelt.org/477b2242-e6a2-4599-8d7d-52d225d39692/c3547552-3022-4c61-8d1c-52b44db4589a", "https://ms.immowelt.org/b6435004-0e35-4fc6-b7d5-548a1c6bd865/18ecd94d-aeb5-4db7-9ece-c9056b7c18f5", "https://ms.immowelt.org/fffca264-e9b1-4647-8e64-c60be93744f8/191162cc-7517-4fc6-be4a-5a945715134c", "https://ms.immowelt.org/93d1fc99-9be1-4826-a8df-e3a22d03c005/b7d69f86-c033-49d5-941d-d0decd68be56", "https://ms.immowelt.org/3593ecfa-644d-4979-ae0d-0fc4d4e9ae7a/1ec892e1-00ff-49d2-b6db-03c288ea7d6d", "https://ms.immowelt.org/19ab9c79-4170-4a92-b48b-ed773d6437fd/7b0670d3-4d77-4ec4-919e-937dab966ff8", "https://ms.immowelt.org/b4d5be21-2515-403f-8215-e2eb79cc92db/3cf8ef03-8de3-42ff-929d-ff3a020a5b25", "https://ms.immowelt.org/b2b6d755-c0f3-4890-a478-6ac51ea07508/68baa9ce-b806-47ca-8c89-b3ed586305cf", "https://ms.immowelt.org/cc54eb68-d41d-4945-85c5-cc8e8c9b5e81/e3d73054-2d37-4c41-8117-5dad34c03029", "https://ms.immowelt.org/8f0aea64-d4d8-47a0-a269-0ddfd36d2981/93df92fb-6eae-41b6-8679-00bb4a30474f", "https://ms.immowelt.org/824566f0-3a4c-4806-a951-d3eecde604c1/9e360642-24e2-44c2-94e2-d1336be6f715", "https://ms.immowelt.org/edc590e2-19f1-4a3d-8938-add9dfc4b56f/db68466a-1681-4bf8-85cb-489b4ffa5436"]'
$ws.Cells.Item(19, 13).Value = '{"companyName": "Michael Engelbrecht", "address": {"city": "W\u00fcrzburg", "zipCode": "97072", "street": "Sartoriusstra\u00dfe 7"}, "salutation": "Herr/Frau", "firstName": "Michael", "lastName": "Engelbrecht", "mobile": "0171-9546505"}'
$ws.Cells.Item(19, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Boutique Wohnung m\u00f6bliert mieten", "basicObjectPricEur": 1100, "basicLivingSpace": 55, "basicConstructionYear": null, "basicCity": "W\u00fcrzburg", "basicUrl": "https://www.immowelt.de/expose/2uewf4g", "basicContactMobile": "0171-9546505"}'

# Row 20 <= original row 21 (url "https://www.immowelt.de/expose/2aebr5t")
$ws.Cells.Item(20, 1).Value = '"https://www.immowelt.de/expose/2aebr5t"'
$ws.Cells.Item(20, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(20, 3).Value = '"Tolle, neuwertige 3,5Zi Erdgeschosswohnung mit Terrasse zwischen Main und Wein in Eibelstadt!"'
$ws.Cells.Item(20, 4).Value = '1290'
$ws.Cells.Item(20, 5).Value = '"EUR"'
$ws.Cells.Item(20, 6).Value = '93'
$ws.Cells.Item(20, 7).Value = '3.5'
$ws.Cells.Item(20, 8).Value = '2022'
$ws.Cells.Item(20, 9).Value = '["Neubau", "barriefrei", "gaestewc", "FENSTER", "Kelleranteil", "Personenaufzug", "PELLET", "Fu\u00dfbodenheizung", "Garten", "Einbauk\u00fcche", "Terrasse"]'
$ws.Cells.Item(20, 10).Value = '{"Country": "Deutschland", "City": "Eibelstadt", "ZipCode": "97246", "Street": "W\u00fcrzburger Strasse 25", "LocationId": 11869, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(20, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "2600", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in Eibelstadt", "Link": "https://www.immowelt.de/immobilienpreise/eibelstadt/mietspiegel"}}, "DataTable": [{"NumberValue": 1290, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 280, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"StringValue": "in Nebenkosten enthalten", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten"}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 1570, "Unit": "EUR"}, {"NumberValue": 40, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "2 Stellpl\u00e4tze, je"}]}'
$ws.Cells.Item(20, 12).Value = '["https://ms.immowelt.org/7c821767-b600-4ad9-8a2a-50874a2a87cf/92a3b2aa-0393-4e55-ad97-c6c3a803a0d8", "https://ms.immowelt.org/8081e4f6-2ad6-46c7-957d-07597db5c5a8/ee5c7a33-80cd-4c89-bb96-98be1b013e27", "https://ms.immowelt.org/515a3a85-6263-4990-b640-f75ec6075c4b/8456cc16-2fc0-409e-ba1f-f103285c6d18", "https://ms.immowelt.org/340c6f6f-746d-421c-aab8-4850f62a5e6b/e971e1c3-7527-4fb5-a83c-f5221d132bc5", "https://ms.immowelt.org/6205eb16-8906-4267-a472-86439d570e9c/39733ff0-9850-43ab-a6fa-f4f2291ba17d", "https://ms.immowelt.org/3b9c4f71-f8fa-46fa-a1d6-7bab696693cb/eccec98b-d098-4d62-ae2c-f62c71c18721", "https://ms.immowelt.org/ab60445b-088f-4fe5-9554-7448a4a29bb2/0eb8642f-4cb3-4252-b4b5-622b757475d6", "https://ms.immowelt.org/6b12d3cb-5e11-4c67-9c0e-a69e674fc018/3d35f5b5-f31d-4200-8572-5eb83a1079e5", "https://ms.immowelt.org/11466902-020f-43ec-9441-09395d524b18/3136e448-448b-41f8-bcec-58fea52abdd5"]'
$ws.Cells.Item(20, 13).Value = '{"salutation": "Herr", "firstName": "Maximilian", "lastName": "Kauschke"}'
$ws.Cells.Item(20, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "Tolle, neuwertige 3,5Zi Erdgeschosswohnung mit Terrasse zwischen Main und Wein in Eibelstadt!", "basicObjectPricEur": 1290, "basicLivingSpace": 93, "basicRooms": 3.5, "basicConstructionYear": 2022, "basicCity": "Eibelstadt", "basicStreet": "W\u00fcrzburger Strasse 25", "basicUrl": "https://www.immowelt.de/expose/2aebr5t"}'

# Row 21 <= original row 18 (url "https://www.immowelt.de/expose/2almm5r")
$ws.Cells.Item(21, 1).Value = '"https://www.immowelt.de/expose/2almm5r"'
$ws.Cells.Item(21, 2).Value = '{"EstateTypeGerman": "WOHNUNG", "DistributionTypeGerman": "ZUR_MIETE", "EstateType": "APARTMENT", "DistributionType": "RENT"}'
$ws.Cells.Item(21, 3).Value = '"TOP-LAGE SANDERAU MIT GROSSER DACHTERRASSE"'
$ws.Cells.Item(21, 4).Value = '1700'
$ws.Cells.Item(21, 5).Value = '"EUR"'
$ws.Cells.Item(21, 6).Value = '153.3'
$ws.Cells.Item(21, 7).Value = '3'
$ws.Cells.Item(21, 8).Value = '1998'
$ws.Cells.Item(21, 9).Value = '["WANNE", "gaestewc", "FENSTER", "Kelleranteil", "Personenaufzug", "Dachgeschoss", "FERNE", "GEPFLEGT", "GAS", "Zentralheizung", "teilweise m\u00f6bliert", "TEPPICH", "PARKETT", "FLIESEN", "Einbauk\u00fcche", "offene K\u00fcche", "Stellplatz", "Tiefgarage", "REINIGUNG", "Terrasse", "Kunststofffenster"]'
$ws.Cells.Item(21, 10).Value = '{"Country": "Deutschland", "City": "W\u00fcrzburg", "District": "Sanderau", "ZipCode": "97072", "Street": "Gartenstra\u00dfe 9", "LocationId": 496022, "PublishStreet": true, "FederalState": "Bayern", "FederalStateId": 2}'
$ws.Cells.Item(21, 11).Value = '{"AdditionalInformation": {"Deposit": {"StringValue": "5.400,00", "Key": "PRICE_DEPOSIT", "Label": "Kaution"}, "MarketPricing": {"Heading": "Mietspiegel in W\u00fcrzburg (Sanderau)", "Link": "https://www.immowelt.de/immobilienpreise/wuerzburg/mietspiegel"}}, "DataTable": [{"NumberValue": 1700, "Unit": "EUR", "Key": "PRICE_RENT_COLD", "Label": "Kaltmiete"}, {"NumberValue": 240, "Unit": "EUR", "Key": "PRICE_ADDITIONALCOSTS", "Label": "Nebenkosten"}, {"NumberValue": 300, "Unit": "EUR", "Key": "PRICE_HEATINGCOSTS", "Label": "Heizkosten", "Comments": ["in Warmmiete enthalten"]}, {"Key": "PRICE_RENT_WARM", "Label": "Warmmiete", "NumberValue": 2340, "Unit": "EUR"}, {"NumberValue": 50, "Unit": "EUR", "Key": "PRICE_PARKINGPRICE", "Label": "2 Stellpl\u00e4tze, je"}]}'
$ws.Cells.Item(21, 12).Value = '["https://ms.immowelt.org/338d1cbc-cb3f-423f-8523-42cc120d72f9/e4b563bc-671f-4dea-a726-71d7b9386009", "https://ms.immowelt.org/4da068e4-1a27-465c-b0e7-ffb77b79880b/908c43fa-aa39-4e85-97c8-98295a3c363c", "https://ms.immowelt.org/cab19ffc-0f4a-42a9-89e9-7189e9dafd92/2a5e6769-ad31-4185-a454-e8d753be9cdb", "https://ms.immowelt.org/f6511334-9f04-476c-9659-9f19a5df115f/ff04d8f5-f770-493d-9394-6ac1a456d4ef", "https://ms.immowelt.org/733f5f85-bc20-4e60-ba25-388eeb6563cf/4094751f-b73c-432d-b6b4-6d68ba45d6d9", "https://ms.immowelt.org/05e1b42e-75a4-4206-aaad-48ffa57f162d/a3aafbff-20b6-4c17-a8a1-b916eb14b8aa", "https://ms.immowelt.org/1233e75d-2ab4-4291-8b67-e87c5a4898f1/6153dd66-9325-4066-80d3-ed8642487cd8", "https://ms.immowelt.org/99f6825d-9789-4ae1-9a78-4e6acb5853a7/f727c1c7-69a0-420d-adf8-9119f5849a16", "https://ms.immowelt.org/19d94014-cd12-4f9a-81a5-8f9745f63b90/48988be8-85bb-4980-a25a-276af2606ac6", "https://ms.immowelt.org/f351ce35-7e50-4939-b806-4f6dc375cc86/84c997d9-1332-4173-9248-fb7b75a10943", "https://ms.immowelt.org/cdd250a7-47f6-4535-ac59-6ea51157c95e/bed7b410-ef26-42d6-a562-595eafba5999", "https://ms.immowelt.org/644f0235-e6ba-482b-bb26-2fea622cd4a1/6d2746f3-ea97-4996-8617-e82aeced34e0", "https://ms.immowelt.org/0174f587-e7b3-4c39-8603-c1d8465305ec/508d427a-f8c5-4d05-bd55-06ee858889be", "https://ms.immowelt.org/e3dbe0a9-b856-43a6-97a0-ca238574c6f8/f20530d2-b140-4eb5-b229-ad329853d7f3"]'
$ws.Cells.Item(21, 13).Value = '{"companyName": "hv Immobilien", "address": {"city": "W\u00fcrzburg", "zipCode": "97082", "street": "Wilhelm-Dahl-Str. 9"}, "profilePicture": "https://filestore.immowelt.de/ProfilBilder/150_98e1685eae1346a2a90a2847489951ff.jpg", "salutation": "Herr", "firstName": "Hartwig", "lastName": "Vogel", "mobile": "0175-5994950", "phone": "0931/45466333"}'
$ws.Cells.Item(21, 14).Value = '{"basicTranactionType": "RENT", "basicTitle": "TOP-LAGE SANDERAU MIT GROSSER DACHTERRASSE", "basicObjectPricEur": 1700, "basicLivingSpace": 153.3, "basicRooms": 3, "basicConstructionYear": 1998, "basicCity": "W\u00fcrzburg", "basicStreet": "Gartenstra\u00dfe 9", "basicUrl": "https://www.immowelt.de/expose/2almm5r", "basicContactPhone": "0931/45466333", "basicContactMobile": "0175-5994950"}'
